# Meta-analysis workbook touch-up: clear stray "NA" placeholders on the
# "Operative time" forest-plot data sheet, and restore "Stroke" as the
# active/selected worksheet tab.

$wb = $excel.ActiveWorkbook

# --- Operative time sheet: clear the cells that only held the "NA" text ---
$ws = $wb.Worksheets.Item("Operative time")

$ws.Range("C2:E2").ClearContents()
$ws.Range("I2").ClearContents()
$ws.Range("K2:M2").ClearContents()
$ws.Range("Q2").ClearContents()

$ws.Range("E3:I3").ClearContents()
$ws.Range("M3:Q3").ClearContents()

$ws.Range("E4:I4").ClearContents()
$ws.Range("M4:Q4").ClearContents()

# Leave the cursor/selection on this sheet at F6
$ws.Activate()
$ws.Range("F6").Select()

# --- Make "Stroke" the active sheet/tab again ---
$stroke = $wb.Worksheets.Item("Stroke")
$stroke.Activate()
